$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Inventory")

$ws.Range('B2').Value = 'Canon imageCLASS 2200 Advanced Copier'
$ws.Range('D2').Value = 7
$ws.Range('E2').Value = 2
$ws.Range('B3').Value = 'Fellowes PB500 Electric Punch Plastic Comb Binding Machine with Manual Bind'
$ws.Range('B4').Value = 'Cisco TelePresence System EX90 Videoconferencing Unit'
$ws.Range('B5').Value = 'HON 5400 Series Task Chairs for Big and Tall'
$ws.Range('B6').Value = 'GBC DocuBind TL300 Electric Binding System'
$ws.Range('B7').Value = 'GBC Ibimaster 500 Manual ProClick Binding System'
$ws.Range('B8').Value = 'HP Designjet T520 Inkjet Large Format Printer - 24" Color'
$ws.Range('B9').Value = 'GBC DocuBind P400 Electric Binding System'
$ws.Range('B10').Value = 'High Speed Automatic Electric Letter Opener'
$ws.Range('B11').Value = 'Lexmark MX611dhe Monochrome Laser Printer'
$ws.Range('B12').Value = 'Hewlett Packard LaserJet 3310 Copier'
$ws.Range('B13').Value = 'Riverside Palais Royal Lawyers Bookcase, Royale Cherry Finish'
$ws.Range('B14').Value = 'Martin Yale Chadless Opener Electric Letter Opener'
$ws.Range('B15').Value = '3D Systems Cube Printer, 2nd Generation, Magenta'
$ws.Range('B16').Value = 'Ibico EPK-21 Electric Binding System'
$ws.Range('B17').Value = 'Apple iPhone 5'
$ws.Range('B18').Value = 'Bretford Rectangular Conference Table Tops'
$ws.Range('B19').Value = 'Samsung Galaxy Mega 6.3'
$ws.Range('B20').Value = 'Canon PC1060 Personal Laser Copier'
$ws.Range('B21').Value = 'Honeywell Enviracaire Portable HEPA Air Cleaner for 17'' x 22'' Room'
$ws.Range('B22').Value = 'Cubify CubeX 3D Printer Double Head Print'
$ws.Range('B23').Value = 'DMI Eclipse Executive Suite Bookcases'
$ws.Range('B24').Value = 'Tennsco 6- and 18-Compartment Lockers'
$ws.Range('B25').Value = 'Plantronics CS510 - Over-the-Head monaural Wireless Headset System'
$ws.Range('B26').Value = 'Global Troy Executive Leather Low-Back Tilter'
$ws.Range('B27').Value = 'Logitech P710e Mobile Speakerphone'
$ws.Range('B28').Value = 'Chromcraft Bull-Nose Wood Oval Conference Tables & Bases'
$ws.Range('B29').Value = 'SAFCO Arco Folding Chair'
$ws.Range('B30').Value = 'Plantronics Savi W720 Multi-Device Wireless Headset System'
$ws.Range('B31').Value = 'Hon Deluxe Fabric Upholstered Stacking Chairs, Rounded Back'
$ws.Range('B32').Value = 'Global Deluxe High-Back Manager''s Chair'
$ws.Range('B33').Value = 'GuestStacker Chair with Chrome Finish Legs'
$ws.Range('B34').Value = 'Hon 4070 Series Pagoda Armless Upholstered Stacking Chairs'
$ws.Range('B35').Value = 'Tennsco Double-Tier Lockers'
$ws.Range('B36').Value = 'Hot File 7-Pocket, Floor Stand'
$ws.Range('B37').Value = 'Tennsco Single-Tier Lockers'
$ws.Range('B38').Value = 'Wilson Electronics DB Pro Signal Booster'
$ws.Range('B39').Value = 'Hewlett Packard 610 Color Digital Copier / Printer'
$ws.Range('B40').Value = 'Okidata MB760 Printer'
$ws.Range('B41').Value = 'Bush Advantage Collection Racetrack Conference Table'
$ws.Range('B42').Value = 'Ativa V4110MDD Micro-Cut Shredder'
$ws.Range('B43').Value = 'Hon 4700 Series Mobuis Mid-Back Task Chairs with Adjustable Arms'
$ws.Range('B44').Value = 'Global Commerce Series High-Back Swivel/Tilt Chairs'
$ws.Range('B45').Value = 'Canon Imageclass D680 Copier / Fax'
